$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel resistor package (was "R-W4")
$ws.Range("C2").Value = "R-1/4W"
$ws.Range("C3").Value = "R-1/4W"
$ws.Range("C4").Value = "R-1/4W"
$ws.Range("C5").Value = "R-1/4W"

# Relabel capacitor packages / descriptions
$ws.Range("C9").Value = "E-P2.5mm 6.3x11.5mm"
$ws.Range("E8").Value = "Film Capacitor THT"
$ws.Range("C8").Value = "C-P5mm 11x7.2mm"

$ws.Range("C6").Value = "C-P5mm"
$ws.Range("C7").Value = "C-P5mm"
$ws.Range("C6").ClearFormats()
$ws.Range("C7").ClearFormats()

$ws.Range("E6").Value = "Ceramic Capacitor THT"
$ws.Range("E7").Value = "Ceramic Capacitor THT"

$ws.Range("E9").Value = "Electrolytic Capacitor THT"

# Update the active selection to match
$ws.Range("D14").Select()
